$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - use Text format to preserve exact formatting
# (trailing zeros, multi-dot "thousands" style numbers, etc.)
$priceUpdates = @{
    2 = "36.408.61"
    3 = "1.978.17"
    5 = "245.59"
    7 = "59.08"
    10 = "56.70"
    11 = "0.0870"
    13 = "22.58"
    14 = "0.858"
    15 = "2.267.71"
    16 = "13.80"
    18 = "1.984.89"
    19 = "36.279.47"
    20 = "0.0₃0905"
    21 = "70.48"
    22 = "5.28"
    23 = "234.45"
    27 = "9.83"
    28 = "164.98"
    30 = "19.88"
    32 = "1.18"
    34 = "0.0655"
    35 = "4.43"
    40 = "2.93"
    41 = "1.21"
    42 = "0.0961"
    45 = "1.08"
    46 = "16.23"
    47 = "91.15"
    48 = "1.365.96"
    49 = "7.44"
    51 = "45.41"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# Volume(1h) (column E) updates - text already keeps padding spaces so it will not be
# misinterpreted as a numeric percentage by Excel
$volumeUpdates = @{
    2 = "  -2.40%  "
    3 = "  -3.07%  "
    4 = "  -0.03%  "
    5 = "  -2.41%  "
    6 = "  -3.82%  "
    7 = "  -9.19%  "
    8 = "  -0.09%  "
    9 = "  -6.78%  "
    10 = "  -4.77%  "
    11 = "  +11.17%  "
    12 = "  +0.04%  "
    13 = "  -2.38%  "
    14 = "  -6.53%  "
    15 = "  -3.21%  "
    16 = "  -6.24%  "
    17 = "  -4.29%  "
    18 = "  -2.87%  "
    19 = "  -2.46%  "
    20 = "  +3.08%  "
    21 = "  -4.03%  "
    22 = "  -3.55%  "
    23 = "  -1.84%  "
    24 = "  -0.01%  "
    25 = "  -3.84%  "
    26 = "  -2.18%  "
    27 = "  -1.45%  "
    28 = "  +2.32%  "
    29 = "  -0.06%  "
    30 = "  -0.44%  "
    31 = "  -1.92%  "
    32 = "  +0.06%  "
    33 = "  -4.96%  "
    34 = "  +4.75%  "
    35 = "  -4.94%  "
    38 = "  -1.77%  "
    39 = "  -6.34%  "
    40 = "  -1.45%  "
    41 = "  -4.27%  "
    42 = "  -4.71%  "
    43 = "  -5.47%  "
    44 = "  -1.96%  "
    45 = "  -5.89%  "
    46 = "  -6.89%  "
    47 = "  -4.66%  "
    48 = "  -2.05%  "
    49 = "  -4.87%  "
    50 = "  -2.34%  "
    51 = "  -3.85%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

# Rows 36 and 37 swapped coins: THORChain <-> BinanceUSD (with refreshed price/volume)
$ws.Cells.Item(36, 2).Value = "BinanceUSD"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$cellD36 = $ws.Cells.Item(36, 4)
$cellD36.NumberFormat = "@"
$cellD36.Value = "1.00"
$ws.Cells.Item(36, 5).Value = "  -0.04%  "

$ws.Cells.Item(37, 2).Value = "THORChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cellD37 = $ws.Cells.Item(37, 4)
$cellD37.NumberFormat = "@"
$cellD37.Value = "6.10"
$ws.Cells.Item(37, 5).Value = "  -3.14%  "

